# Fix summary of capital gains
#
# - "Shares" sheet: rename "Total Gain [EUR]" header to "Gain [EUR]" and
#   drop the trailing "Total" sum row (it duplicated what the ELSTER summary
#   already computes).
# - Same "Total" row removal on "Foreign Currencies", "Dividend Payments",
#   "Fees" and "Tax Withholdings".
# - "ELSTER - Summary": tidy up the duplicated wording in the "Zeile 23"
#   description, and fix the gain/loss split so "Zeile 20" reflects only the
#   sum of the winning trades and "Zeile 23" only the losing ones (instead of
#   the previous net total / 0 split).

$wb = $excel.ActiveWorkbook

$shares = $wb.Worksheets.Item("Shares")
$shares.Range("I1").Value = "Gain [EUR]"
$shares.Rows.Item(11).Delete()

$wb.Worksheets.Item("Foreign Currencies").Rows.Item(12).Delete()
$wb.Worksheets.Item("Dividend Payments").Rows.Item(4).Delete()
$wb.Worksheets.Item("Fees").Rows.Item(11).Delete()
$wb.Worksheets.Item("Tax Withholdings").Rows.Item(4).Delete()

$summary = $wb.Worksheets.Item("ELSTER - Summary")
$summary.Range("B4").Value = "Zeile 23: In den Zeilen 18 und 19 enthaltene Verluste aus der Veräuerung von Aktien i. S. d. § 20 Abs. 2 Satz 1 Nr. 1 EStG"
$summary.Range("C3").Value = 974.86
$summary.Range("C4").Value = 247.01
